$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025624586480298
$ws.Range("D2").Value = 1.033405549184609
$ws.Range("E2").Value = 1.025954472687689
$ws.Range("F2").Value = 1.037173900432384
$ws.Range("I2").Value = 1.035929501502842
$ws.Range("J2").Value = 1.030792566634761
$ws.Range("K2").Value = 1.036207992280795
$ws.Range("L2").Value = 1.028778532994936
$ws.Range("M2").Value = 1.039965539177091
$ws.Range("N2").Value = 1.032256410312949
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026557352035538
$ws.Range("D3").Value = 1.033926613548358
$ws.Range("E3").Value = 1.026745781259389
$ws.Range("F3").Value = 1.038345384730879
$ws.Range("I3").Value = 1.036148148833666
$ws.Range("J3").Value = 1.031364831758593
$ws.Range("K3").Value = 1.036538799109067
$ws.Range("L3").Value = 1.02937729738806
$ws.Range("M3").Value = 1.04094581874574
$ws.Range("N3").Value = 1.03282948811889
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027161159777449
$ws.Range("D4").Value = 1.034263762112921
$ws.Range("E4").Value = 1.027258408185518
$ws.Range("F4").Value = 1.039103866130239
$ws.Range("I4").Value = 1.036288288934316
$ws.Range("J4").Value = 1.031734770078981
$ws.Range("K4").Value = 1.036752125522592
$ws.Range("L4").Value = 1.029764682354929
$ws.Range("M4").Value = 1.041580006727537
$ws.Range("N4").Value = 1.033199951794125
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027415058809269
$ws.Range("D5").Value = 1.034405494291105
$ws.Range("E5").Value = 1.027474058468713
$ws.Range("F5").Value = 1.039422840144032
$ws.Range("I5").Value = 1.036346882638335
$ws.Range("J5").Value = 1.031890206374364
$ws.Range("K5").Value = 1.036841632783871
$ws.Range("L5").Value = 1.029927524896414
$ws.Range("M5").Value = 1.041846590893657
$ws.Range("N5").Value = 1.033355608826877
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027457692971336
$ws.Range("D6").Value = 1.034429291405752
$ws.Range("E6").Value = 1.027510275390987
$ws.Range("F6").Value = 1.039476403644662
$ws.Range("I6").Value = 1.036356701927993
$ws.Range("J6").Value = 1.031916299762925
$ws.Range("K6").Value = 1.036856651147732
$ws.Range("L6").Value = 1.02995486602624
$ws.Range("M6").Value = 1.041891349876214
$ws.Range("N6").Value = 1.033381739271042
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027164552160635
$ws.Range("D7").Value = 1.034265655967084
$ws.Range("E7").Value = 1.027261289158862
$ws.Range("F7").Value = 1.039108127850485
$ws.Range("I7").Value = 1.036289073129025
$ws.Range("J7").Value = 1.031736847364247
$ws.Range("K7").Value = 1.036753322212461
$ws.Range("L7").Value = 1.029766858321317
$ws.Range("M7").Value = 1.041583568950176
$ws.Range("N7").Value = 1.033202032029374
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025939767628826
$ws.Range("D8").Value = 1.033581647511916
$ws.Range("E8").Value = 1.026221774782424
$ws.Range("F8").Value = 1.037569715326242
$ws.Range("I8").Value = 1.036003671379615
$ws.Range("J8").Value = 1.03098603951383
$ws.Range("K8").Value = 1.03631994011082
$ws.Range("L8").Value = 1.028980899275884
$ws.Range("M8").Value = 1.040296853709084
$ws.Range("N8").Value = 1.032450157945691
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.023783453707819
$ws.Range("D9").Value = 1.032376289203077
$ws.Range("E9").Value = 1.024394640564982
$ws.Range("F9").Value = 1.034862300497709
$ws.Range("I9").Value = 1.035490523547371
$ws.Range("J9").Value = 1.02966032703953
$ws.Range("K9").Value = 1.035550733481907
$ws.Range("L9").Value = 1.027595547167567
$ws.Range("M9").Value = 1.038028594069461
$ws.Range("N9").Value = 1.031122562807616
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.022347233599032
$ws.Range("D10").Value = 1.031572773646327
$ws.Range("E10").Value = 1.023179722622176
$ws.Range("F10").Value = 1.033059679642205
$ws.Range("I10").Value = 1.035141576657848
$ws.Range("J10").Value = 1.028774746505995
$ws.Range("K10").Value = 1.035034265113197
$ws.Range("L10").Value = 1.02667176070826
$ws.Range("M10").Value = 1.036515816633923
$ws.Range("N10").Value = 1.030235724648146
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021725655151515
$ws.Range("D11").Value = 1.031224874615155
$ws.Range("E11").Value = 1.022654414988688
$ws.Range("F11").Value = 1.03227967315954
$ws.Range("I11").Value = 1.034988861007674
$ws.Range("J11").Value = 1.028390867726937
$ws.Range("K11").Value = 1.034809771095088
$ws.Range("L11").Value = 1.026271708326757
$ws.Range("M11").Value = 1.035860623601897
$ws.Range("N11").Value = 1.029851300717191
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021494820955207
$ws.Range("D12").Value = 1.031095655200043
$ws.Range("E12").Value = 1.022459407514664
$ws.Range("F12").Value = 1.031990024847864
$ws.Range("I12").Value = 1.034931892812482
$ws.Range("J12").Value = 1.028248216041158
$ws.Range("K12").Value = 1.034726255748156
$ws.Range("L12").Value = 1.026123104730511
$ws.Range("M12").Value = 1.035617232991665
$ws.Range("N12").Value = 1.029708446449655
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021544333506826
$ws.Range("D13").Value = 1.031123372918004
$ws.Range("E13").Value = 1.0225012320519
$ws.Range("F13").Value = 1.032052151741838
$ws.Range("I13").Value = 1.034944123663301
$ws.Range("J13").Value = 1.028278818109182
$ws.Range("K13").Value = 1.034744175863736
$ws.Range("L13").Value = 1.026154980981643
$ws.Range("M13").Value = 1.03566944213539
$ws.Range("N13").Value = 1.029739091976126
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021706573348632
$ws.Range("D14").Value = 1.031214193171574
$ws.Range("E14").Value = 1.022638293254773
$ws.Range("F14").Value = 1.032255729092872
$ws.Range("I14").Value = 1.034984156949193
$ws.Range("J14").Value = 1.028379077355889
$ws.Range("K14").Value = 1.034802870301979
$ws.Range("L14").Value = 1.026259424832244
$ws.Range("M14").Value = 1.035840505322974
$ws.Range("N14").Value = 1.029839493602464
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021806541007122
$ws.Range("D15").Value = 1.031270151337035
$ws.Range("E15").Value = 1.022722756476852
$ws.Range("F15").Value = 1.032381170517845
$ws.Range("I15").Value = 1.035008790618638
$ws.Range("J15").Value = 1.028440842186178
$ws.Range("K15").Value = 1.034839016912073
$ws.Range("L15").Value = 1.026323775317713
$ws.Range("M15").Value = 1.035945899980278
$ws.Range("N15").Value = 1.029901346145896
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.022388492377911
$ws.Range("D16").Value = 1.031595863287422
$ws.Range("E16").Value = 1.023214601655333
$ws.Range("F16").Value = 1.033111457456805
$ws.Range("I16").Value = 1.035151677816291
$ws.Range("J16").Value = 1.028800214548341
$ws.Range("K16").Value = 1.035049146003874
$ws.Range("L16").Value = 1.026698309941608
$ws.Range("M16").Value = 1.036559296455871
$ws.Range("N16").Value = 1.030261228858034
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022753619599708
$ws.Range("D17").Value = 1.031800182476823
$ws.Range("E17").Value = 1.023523327290147
$ws.Range("F17").Value = 1.033569691557634
$ws.Range("I17").Value = 1.035240873978587
$ws.Range("J17").Value = 1.029025528212546
$ws.Range("K17").Value = 1.035180724875871
$ws.Range("L17").Value = 1.026933233590566
$ws.Range("M17").Value = 1.036944023483488
$ws.Range("N17").Value = 1.030486862493486
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022966622426796
$ws.Range("D18").Value = 1.031919361054064
$ws.Range("E18").Value = 1.023703474976835
$ws.Range("F18").Value = 1.033837023904197
$ws.Range("I18").Value = 1.03529274431523
$ws.Range("J18").Value = 1.029156909630357
$ws.Range("K18").Value = 1.035257389586573
$ws.Range("L18").Value = 1.027070256023383
$ws.Range("M18").Value = 1.03716841374603
$ws.Range("N18").Value = 1.030618430487979
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023039256010229
$ws.Range("D19").Value = 1.031959998273174
$ws.Range("E19").Value = 1.02376491307575
$ws.Range("F19").Value = 1.033928186181516
$ws.Range("I19").Value = 1.035310404237732
$ws.Range("J19").Value = 1.029201700457516
$ws.Range("K19").Value = 1.0352835161627
$ws.Range("L19").Value = 1.027116976337912
$ws.Range("M19").Value = 1.037244922571235
$ws.Range("N19").Value = 1.030663284923251
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022714441767066
$ws.Range("D20").Value = 1.031778260676466
$ws.Range("E20").Value = 1.023490196367641
$ws.Range("F20").Value = 1.033520521983667
$ws.Range("I20").Value = 1.03523132023532
$ws.Range("J20").Value = 1.029001358345796
$ws.Range("K20").Value = 1.035166616288675
$ws.Range("L20").Value = 1.026908028972626
$ws.Range("M20").Value = 1.036902747415554
$ws.Range("N20").Value = 1.030462658302752
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.021658796437837
$ws.Range("D21").Value = 1.031187448699159
$ws.Range("E21").Value = 1.022597928959312
$ws.Range("F21").Value = 1.032195778408548
$ws.Range("I21").Value = 1.03497237484296
$ws.Range("J21").Value = 1.028349555213981
$ws.Range("K21").Value = 1.034785589786565
$ws.Range("L21").Value = 1.026228668890591
$ws.Range("M21").Value = 1.035790132106926
$ws.Range("N21").Value = 1.029809929535728
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020995346628376
$ws.Range("D22").Value = 1.030816015907011
$ws.Range("E22").Value = 1.022037591864936
$ws.Range("F22").Value = 1.031363327011287
$ws.Range("I22").Value = 1.034808161053469
$ws.Range("J22").Value = 1.027939382196503
$ws.Range("K22").Value = 1.034545281598072
$ws.Range("L22").Value = 1.025801491965961
$ws.Range("M22").Value = 1.035090455032677
$ws.Range("N22").Value = 1.029399174025518
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.021347027570571
$ws.Range("D23").Value = 1.031012915724757
$ws.Range("E23").Value = 1.022334573614916
$ws.Range("F23").Value = 1.031804581049609
$ws.Range("I23").Value = 1.034895346862991
$ws.Range("J23").Value = 1.028156856467362
$ws.Range("K23").Value = 1.034672743519936
$ws.Range("L23").Value = 1.026027949743405
$ws.Range("M23").Value = 1.035461379617285
$ws.Range("N23").Value = 1.029616957134785
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022732144449604
$ws.Range("D24").Value = 1.031788166186414
$ws.Range("E24").Value = 1.023505166579172
$ws.Range("F24").Value = 1.033542739436313
$ws.Range("I24").Value = 1.035235637643556
$ws.Range("J24").Value = 1.029012279792515
$ws.Range("K24").Value = 1.035172991608535
$ws.Range("L24").Value = 1.026919417868667
$ws.Range("M24").Value = 1.036921398338583
$ws.Range("N24").Value = 1.030473595259179
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024340682110795
$ws.Range("D25").Value = 1.032687899717611
$ws.Range("E25").Value = 1.024866444180136
$ws.Range("F25").Value = 1.035561822288039
$ws.Range("I25").Value = 1.035624393973733
$ws.Range("J25").Value = 1.030003370670508
$ws.Range("K25").Value = 1.035750241675002
$ws.Range("L25").Value = 1.027953735680976
$ws.Range("M25").Value = 1.038615100284672
$ws.Range("N25").Value = 1.031466093599898
